$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = "ML195059"
$ws.Range("B19").Value = "Oldie"
$ws.Range("A19").Value = "Oil Room"

$ws.Range("C20").Value = "MJ190517"
$ws.Range("A20").Value = "Boot Exchange"
$ws.Range("B20").Value = "RTE 1"
